# Apply "all functions ok with colons checking" edit:
# - Set C4 to "x" (was a formula-like note "+check all arg_check")
# - Set C10 to "x" (was "issue 82")
# - Fill C5:C8, C12:C17, C19:C22 with "x"
# - Move the active selection to D24

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells that change content
$ws.Range("C4").Value = "x"
$ws.Range("C10").Value = "x"

# Fill in the newly populated cells in column C
$rows = @(5,6,7,8,12,13,14,15,16,17,19,20,21,22)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "x"
}

# Move selection to D24 as recorded in the saved view state
$ws.Range("D24").Select()
